$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C2:C13) from 45185 (2023-09-16) to 45204 (2023-10-05)
$ws.Range("C2:C13").Value = 45204
